$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: add rows 25-27 ---
$ws1.Range("A25").Value = "adding features to biquad, leveller"
$ws1.Range("J25").Value = 20000
$ws1.Range("K25").Formula = "=J24-J25"

$ws1.Range("J26").Value = 19920
$ws1.Range("K26").Formula = "=J25-J26"

$ws1.Range("A27").Value = "optimizing biquad"
$ws1.Range("J27").Value = 19324

# --- Add Sheet2 after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("B1").Value = "size"
$ws2.Range("A1").Value = "what"

$ws2.Range("A2").Value = "with lp/hp sepa"
$ws2.Range("B2").Value = 21899

$ws2.Range("B3").Value = 21895

$ws2.Range("B4").Value = 21899

$ws2.Range("A5").Value = "after biquad optimizations and improvements"
$ws2.Range("B5").Value = 19324

$ws2.Range("A6").Value = "no smasher"
$ws2.Range("B6").Value = 17159

$ws2.Range("A7").Value = "with maj7comp"
$ws2.Range("B7").Value = 19868
$ws2.Range("C7").Formula = "=B7-B6"

$ws2.Range("A8").Value = "removing features"
$ws2.Range("B8").Value = 19872
$ws2.Range("D8").Value = "wtf how did REMOVING features increase size?"

$ws2.Range("A9").Value = "removing a call to exp2"
$ws2.Range("B9").Value = 19868

$ws2.Columns.Item(1).ColumnWidth = 42.7109375

# --- View settings ---
$ws1.Range("J15").Select() | Out-Null
$ws2.Range("A10").Select() | Out-Null
